$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductLoanInput")

# Row 17 (repaymentstrategy) changes from "Mifos style" to the new
# "Penalties, Fees, Interest, Principal order" scenario value, and picks
# up the same left/top-aligned, non-bold formatting used by B1.
$ws.Range("B17").Value = "Penalties, Fees, Interest, Principal order"

$ws.Range("B1").Copy()
$ws.Range("B17").PasteSpecial(-4122)

# Move the active selection to the newly edited cell.
$ws.Range("B17").Select()
